$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": a new client "IMPORTELECTRIC S.A.S" is inserted
# (alphabetically) right before "INTERNEGOCIOS DE HIERRO S.A." at row 33,
# pushing every following client row down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(33).Insert()

$ws1.Range("A33").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Range("B33").Value = "IMPORTELECTRIC S.A.S"
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($c in $cols1) {
    $ws1.Range($c + "33").Value = 0
}

# Genuine data correction: BORJA TORRES LETTY JANET's PORCELANATO total drops.
$ws1.Range("M11").Value = 15614.03

# The final "n de 58" tally row (old row 60, now row 61) must reflect the new
# client count of 59.
$tallyRow1 = 61
$tallyValues1 = @{
    "C" = "0 de 59"; "D" = "2 de 59"; "E" = "1 de 59"; "F" = "0 de 59";
    "G" = "0 de 59"; "H" = "0 de 59"; "I" = "0 de 59"; "J" = "0 de 59";
    "K" = "0 de 59"; "L" = "2 de 59"; "M" = "5 de 59"; "N" = "0 de 59";
    "O" = "0 de 59"; "P" = "1 de 59"; "Q" = "0 de 59"; "R" = "0 de 59";
}
foreach ($c in $cols1) {
    $ws1.Range($c + $tallyRow1).Value = $tallyValues1[$c]
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same new client inserted at the same position.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(33).Insert()

$ws2.Range("A33").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Range("B33").Value = "IMPORTELECTRIC S.A.S"
$cols2 = @("C","D","E","F","G")
foreach ($c in $cols2) {
    $ws2.Range($c + "33").Value = 0
}

# Matching "noviembre" correction for BORJA TORRES LETTY JANET.
$ws2.Range("F11").Value = 23858.61

# Recompute the monthly totals row (old row 60, now row 61); only the
# "noviembre" column actually changes because of the F11 correction above.
$ws2.Range("C61").Value = 88332.23999999999
$ws2.Range("D61").Value = 69222.88
$ws2.Range("E61").Value = 54885.7
$ws2.Range("F61").Value = 33920.9
$ws2.Range("G61").Value = 52000

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO and TOTAL rows shift their
# VENTA / POR CUMPLIR / CUMPLIMIENTO figures because of the PORCELANATO sale
# correction above.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 21058.27
$ws3.Range("E12").Value = 26982.73
$ws3.Range("F12").Value = 0.4383395433067588
$ws3.Range("D14").Value = 34032.78
$ws3.Range("E14").Value = 23854.57196497848
$ws3.Range("F14").Value = 0.5879139197901061
